$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 63.211268
$ws.Range("H2").Value2 = 189.633804
$ws.Range("I2").Value2 = 0.4922609885657722
$ws.Range("J2").Value2 = 0.4922609885657722
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 2.172943
$ws.Range("N2").Value2 = 6.518829
$ws.Range("O2").Value2 = 0.1731236386970244
$ws.Range("P2").Value2 = 0.1731236386970244
$ws.Range("Q2").Value2 = 137.354482321724
$ws.Range("R2").Value2 = 1236.190340895516
$ws.Range("S2").Value2 = 0.0852220135291008
$ws.Range("T2").Value2 = 0.08522201352910082

# Row 3
$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 63.211268
$ws.Range("H3").Value2 = 189.633804
$ws.Range("I3").Value2 = 0.4922609885657722
$ws.Range("J3").Value2 = 0.4922609885657722
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 8.934372
$ws.Range("N3").Value2 = 26.803116
$ws.Range("O3").Value2 = 0.7118230851489483
$ws.Range("P3").Value2 = 0.7118230851489484
$ws.Range("Q3").Value2 = 564.752982903696
$ws.Range("R3").Value2 = 5082.776846133263
$ws.Range("S3").Value2 = 0.3504027355793591
$ws.Range("T3").Value2 = 0.3504027355793592

# Row 4
$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 63.211268
$ws.Range("H4").Value2 = 189.633804
$ws.Range("I4").Value2 = 0.4922609885657722
$ws.Range("J4").Value2 = 0.4922609885657722
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 1.444079
$ws.Range("N4").Value2 = 4.332237
$ws.Range("O4").Value2 = 0.1150532761540272
$ws.Range("P4").Value2 = 0.1150532761540272
$ws.Range("Q4").Value2 = 91.282064682172
$ws.Range("R4").Value2 = 821.538582139548
$ws.Range("S4").Value2 = 0.05663623945731221
$ws.Range("T4").Value2 = 0.05663623945731221

# Row 5
$ws.Range("E5").Value2 = 3
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 43.30706799999999
$ws.Range("H5").Value2 = 129.921204
$ws.Range("I5").Value2 = 0.3372560111523963
$ws.Range("J5").Value2 = 0.3372560111523963
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 2.172943
$ws.Range("N5").Value2 = 6.518829
$ws.Range("O5").Value2 = 0.1731236386970244
$ws.Range("P5").Value2 = 0.1731236386970244
$ws.Range("Q5").Value2 = 94.10379026112399
$ws.Range("R5").Value2 = 846.9341123501159
$ws.Range("S5").Value2 = 0.0583869878231471
$ws.Range("T5").Value2 = 0.05838698782314711

# Row 6
$ws.Range("E6").Value2 = 3
$ws.Range("F6").Value2 = 1
$ws.Range("G6").Value2 = 43.30706799999999
$ws.Range("H6").Value2 = 129.921204
$ws.Range("I6").Value2 = 0.3372560111523963
$ws.Range("J6").Value2 = 0.3372560111523963
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 8.934372
$ws.Range("N6").Value2 = 26.803116
$ws.Range("O6").Value2 = 0.7118230851489483
$ws.Range("P6").Value2 = 0.7118230851489484
$ws.Range("Q6").Value2 = 386.9214557412959
$ws.Range("R6").Value2 = 3482.293101671663
$ws.Range("S6").Value2 = 0.2400666143435269
$ws.Range("T6").Value2 = 0.2400666143435269

# Row 7
$ws.Range("E7").Value2 = 3
$ws.Range("F7").Value2 = 1
$ws.Range("G7").Value2 = 43.30706799999999
$ws.Range("H7").Value2 = 129.921204
$ws.Range("I7").Value2 = 0.3372560111523963
$ws.Range("J7").Value2 = 0.3372560111523963
$ws.Range("K7").Value2 = 3
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 1.444079
$ws.Range("N7").Value2 = 4.332237
$ws.Range("O7").Value2 = 0.1150532761540272
$ws.Range("P7").Value2 = 0.1150532761540272
$ws.Range("Q7").Value2 = 62.538827450372
$ws.Range("R7").Value2 = 562.849447053348
$ws.Range("S7").Value2 = 0.03880240898572233
$ws.Range("T7").Value2 = 0.03880240898572233

# Row 8
$ws.Range("E8").Value2 = 3
$ws.Range("F8").Value2 = 1
$ws.Range("G8").Value2 = 21.891734
$ws.Range("H8").Value2 = 65.675202
$ws.Range("I8").Value2 = 0.1704830002818315
$ws.Range("J8").Value2 = 0.1704830002818315
$ws.Range("K8").Value2 = 3
$ws.Range("L8").Value2 = 1
$ws.Range("M8").Value2 = 2.172943
$ws.Range("N8").Value2 = 6.518829
$ws.Range("O8").Value2 = 0.1731236386970244
$ws.Range("P8").Value2 = 0.1731236386970244
$ws.Range("Q8").Value2 = 47.569490153162
$ws.Range("R8").Value2 = 428.125411378458
$ws.Range("S8").Value2 = 0.02951463734477651
$ws.Range("T8").Value2 = 0.02951463734477651

# Row 9
$ws.Range("E9").Value2 = 3
$ws.Range("F9").Value2 = 1
$ws.Range("G9").Value2 = 21.891734
$ws.Range("H9").Value2 = 65.675202
$ws.Range("I9").Value2 = 0.1704830002818315
$ws.Range("J9").Value2 = 0.1704830002818315
$ws.Range("K9").Value2 = 3
$ws.Range("L9").Value2 = 1
$ws.Range("M9").Value2 = 8.934372
$ws.Range("N9").Value2 = 26.803116
$ws.Range("O9").Value2 = 0.7118230851489483
$ws.Range("P9").Value2 = 0.7118230851489484
$ws.Range("Q9").Value2 = 195.588895281048
$ws.Range("R9").Value2 = 1760.300057529432
$ws.Range("S9").Value2 = 0.1213537352260623
$ws.Range("T9").Value2 = 0.1213537352260623

# Row 10
$ws.Range("E10").Value2 = 3
$ws.Range("F10").Value2 = 1
$ws.Range("G10").Value2 = 21.891734
$ws.Range("H10").Value2 = 65.675202
$ws.Range("I10").Value2 = 0.1704830002818315
$ws.Range("J10").Value2 = 0.1704830002818315
$ws.Range("K10").Value2 = 3
$ws.Range("L10").Value2 = 1
$ws.Range("M10").Value2 = 1.444079
$ws.Range("N10").Value2 = 4.332237
$ws.Range("O10").Value2 = 0.1150532761540272
$ws.Range("P10").Value2 = 0.1150532761540272
$ws.Range("Q10").Value2 = 31.613393342986
$ws.Range("R10").Value2 = 284.520540086874
$ws.Range("S10").Value2 = 0.01961462771099266
$ws.Range("T10").Value2 = 0.01961462771099265
